$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.127483129501343
$ws.Range("B1").Value = 3.747226238250732
$ws.Range("C1").Value = 3.146496772766113
$ws.Range("D1").Value = 2.04565954208374
$ws.Range("E1").Value = 1.170873522758484
